# Auto-generated edit script applying numeric corrections to Profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("J76").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("J79").Value = 0
$ws.Range("I106").Value = 3634
$ws.Range("K106").Value = 3634
$ws.Range("M106").Value = -3003
$ws.Range("H106").Value = 3634
$ws.Range("I113").Value = 4117
$ws.Range("K113").Value = 4117
$ws.Range("M113").Value = -863
$ws.Range("H113").Value = 4362.1665
$ws.Range("M125").Value = -42090
$ws.Range("K125").Value = 44550
$ws.Range("I125").Value = 4950
$ws.Range("H125").Value = 8474.75
$ws.Range("H129").Value = 2408.0356
$ws.Range("K129").Value = 2482.38465
$ws.Range("I129").Value = 827.46155
$ws.Range("M129").Value = 2517.61535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -1751.0834
$ws.Range("H2").Value = 1738.1333
$ws.Range("I2").Value = 1864.0834
$ws.Range("K2").Value = 1864.0834
$ws.Range("H32").Value = 8862.158
$ws.Range("M32").Value = -7803.647
$ws.Range("K32").Value = 8090.647
$ws.Range("I32").Value = 8090.647
$ws.Range("N88").Value = -3140.8
$ws.Range("L88").Value = 2328.8
$ws.Range("H88").Value = 1767.4546
$ws.Range("J88").Value = 2328.8
$ws.Range("H91").Value = 1767.4546
$ws.Range("J91").Value = 2328.8
$ws.Range("L91").Value = 2328.8
$ws.Range("N91").Value = -5136.8
$ws.Range("H110").Value = 3180.8
$ws.Range("K110").Value = 452
$ws.Range("I110").Value = 452
$ws.Range("M110").Value = 1593
$ws.Range("K116").Value = 1864.0834
$ws.Range("H116").Value = 1738.1333
$ws.Range("M116").Value = 429.9166
$ws.Range("I116").Value = 1864.0834
$ws.Range("H122").Value = 2997
$ws.Range("M122").Value = -6541
$ws.Range("I122").Value = 2997
$ws.Range("K122").Value = 8991

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1738.1333
$ws.Range("K3").Value = 1864.0834
$ws.Range("I3").Value = 1864.0834
$ws.Range("M3").Value = -1750.0834
$ws.Range("M75").Value = -5210.857
$ws.Range("K75").Value = 6146.857
$ws.Range("I75").Value = 6146.857
$ws.Range("H75").Value = 21628.5
$ws.Range("H78").Value = 21628.5
$ws.Range("I78").Value = 6146.857
$ws.Range("M78").Value = -13760.571
$ws.Range("K78").Value = 18440.571
$ws.Range("H94").Value = 1278.8889
$ws.Range("J94").Value = 1099
$ws.Range("N94").Value = -2001
$ws.Range("L94").Value = 1099

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2074.75
$ws.Range("M122").Value = -3548.5
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("K134").Value = 26598.375
$ws.Range("H134").Value = 8972.9
$ws.Range("L134").Value = 28200
$ws.Range("J134").Value = 9400
$ws.Range("I134").Value = 8866.125
$ws.Range("M134").Value = -24063.375
$ws.Range("N134").Value = -33270
$ws.Range("L141").Value = 94093.14
$ws.Range("N141").Value = -104453.14
$ws.Range("H141").Value = 94093.14
$ws.Range("J141").Value = 94093.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 2900
$ws.Range("K110").Value = 8700
$ws.Range("I110").Value = 2900
$ws.Range("M110").Value = -4610
$ws.Range("H136").Value = 1602.3334
$ws.Range("K136").Value = 4807.0002
$ws.Range("M136").Value = 292.9997999999996
$ws.Range("I136").Value = 1602.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("J62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N63").Value = -16371
$ws.Range("L63").Value = 14999
$ws.Range("H63").Value = 14999
$ws.Range("J63").Value = 14999
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("H65").Value = 30000
$ws.Range("N65").Value = -96864
$ws.Range("N66").Value = -51861
$ws.Range("H66").Value = 14999
$ws.Range("L66").Value = 44997
$ws.Range("J66").Value = 14999
$ws.Range("L80").Value = 1703.35
$ws.Range("J80").Value = 1703.35
$ws.Range("M80").Value = -728.4000000000001
$ws.Range("H80").Value = 1707.96
$ws.Range("K80").Value = 1726.4
$ws.Range("N80").Value = -3699.35
$ws.Range("I80").Value = 1726.4
$ws.Range("M83").Value = -3640
$ws.Range("I83").Value = 1726.4
$ws.Range("J83").Value = 1703.35
$ws.Range("N83").Value = -18500.75
$ws.Range("L83").Value = 8516.75
$ws.Range("H83").Value = 1707.96
$ws.Range("K83").Value = 8632
$ws.Range("H87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("N90").Value = -87480
$ws.Range("H90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("M102").Value = 998.4545
$ws.Range("I102").Value = 623.5455
$ws.Range("K102").Value = 623.5455
$ws.Range("H102").Value = 989.9286
$ws.Range("H122").Value = 11367213
$ws.Range("M122").Value = -41674666
$ws.Range("I122").Value = 13892372
$ws.Range("K122").Value = 41677116

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("N7").ClearContents()
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("I40").Value = 9000
$ws.Range("H40").Value = 9836.833
$ws.Range("M40").Value = -8864
$ws.Range("K40").Value = 9000
$ws.Range("N56").Value = -61382
$ws.Range("H56").Value = 60000
$ws.Range("K56").Value = 0
$ws.Range("J56").Value = 60000
$ws.Range("I56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("L56").Value = 60000
$ws.Range("H74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("M82").Value = -331.4
$ws.Range("I82").Value = 692.4
$ws.Range("N82").Value = -2227.5555
$ws.Range("J82").Value = 1505.5555
$ws.Range("H82").Value = 1215.1428
$ws.Range("K82").Value = 692.4
$ws.Range("L82").Value = 1505.5555
$ws.Range("M85").Value = 555.6
$ws.Range("L85").Value = 1505.5555
$ws.Range("H85").Value = 1215.1428
$ws.Range("I85").Value = 692.4
$ws.Range("K85").Value = 692.4
$ws.Range("J85").Value = 1505.5555
$ws.Range("N85").Value = -4001.5555
$ws.Range("H122").Value = 4126.25
$ws.Range("M122").Value = -7756.599999999999
$ws.Range("I122").Value = 3402.2
$ws.Range("K122").Value = 10206.6
$ws.Range("N126").ClearContents()
$ws.Range("L126").Value = 0
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("H136").Value = 3503.3333
$ws.Range("K136").Value = 10509.9999
$ws.Range("M136").Value = -7959.999899999999
$ws.Range("I136").Value = 3503.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M58").Value = -5692
$ws.Range("K58").Value = 6000
$ws.Range("I58").Value = 6000
$ws.Range("H58").Value = 6000
$ws.Range("J70").Value = 50000
$ws.Range("N70").Value = -50630
$ws.Range("L70").Value = 50000
$ws.Range("H70").Value = 50000
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184
$ws.Range("L75").Value = 32999
$ws.Range("N75").Value = -34871
$ws.Range("H75").Value = 32999
$ws.Range("J75").Value = 32999
$ws.Range("H78").Value = 32999
$ws.Range("J78").Value = 32999
$ws.Range("L78").Value = 98997
$ws.Range("N78").Value = -108357
$ws.Range("H103").Value = 16834
$ws.Range("N103").Value = -19178
$ws.Range("J103").Value = 16834
$ws.Range("L103").Value = 16834
$ws.Range("L119").Value = 78000
$ws.Range("J119").Value = 78000
$ws.Range("N119").Value = -87676
$ws.Range("H119").Value = 78000

